$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A95").Value = "EmailCannotBeEmpty"
$ws.Range("B95").Value = "Lutfen E-Posta Giriniz"

$ws.Range("A96").Value = "PasswordCannotBeEmpty"
$ws.Range("B96").Value = "Lutfen Parola Giriniz"

$ws.Range("A97").Value = "invalid_grant"
$ws.Range("B97").Value = "Gecersiz kullanici adi veya parola"
